$d = $word.ActiveDocument

# The "Info:" paragraph for the ProLo Systems project (matthewia/projects/5)
# ends with a sentence describing the UI design / programming work. Expand
# that sentence with more detail about the technologies used and the
# README / project leadership work, per the commit.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$old = " I designed the UI and worked collaboratively to program the application. "
$new = " I designed the UI, and worked collaboratively to program the application using JavaScript and JQuery. Working on this project was my first experience leading a collaborative programming effort. I set up the project and wrote an extensive README with the goal of getting us started by translating my design work into code, and outlining features and specs.  "

$found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target sentence to replace"
}
